# New Shunt Cal sheets for the small scale LBCBs
#
# The "Previous Slopes:" rows on the "Displacement Cals" sheet (rows 20, 41
# and 62 -- one per calibration table) hold hard-coded slope values that are
# compared against the freshly computed "New Slope:" values a few rows above
# (rows 17, 38, 59). Update the four numeric cells (C/D/F/G) in each of those
# rows to the latest shunt-cal results; the "Delta" and "% Difference" rows
# just below each (21/22, 42/43, 63/64) are formulas that recalculate
# automatically off of these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Displacement Cals")

# --- Table 1 (rows 2-22): "Previous Slopes:" row 20 ---
$ws.Range("C20").Value = -0.27018038721863691
$ws.Range("D20").Value = 0.26969468103565636
$ws.Range("F20").Value = -0.27215909090909091
$ws.Range("G20").Value = 0.27176828870122743

# --- Table 2 (rows 23-43): "Previous Slopes:" row 41 ---
$ws.Range("C41").Value = -0.13653794642857142
$ws.Range("D41").Value = 0.1362724941836006
$ws.Range("F41").Value = -0.13740625000000001
$ws.Range("G41").Value = 0.13849747226772158

# --- Table 3 (rows 44-64): "Previous Slopes:" row 62 ---
$ws.Range("C62").Value = -0.13370870535714285
$ws.Range("D62").Value = 0.13209961629671843
$ws.Range("F62").Value = -0.13216517857142857
$ws.Range("G62").Value = 0.13172485577531828

$wb.Save()
